$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the summary table (column M, rows 19-24) with new attendance counts
$ws.Range("M19").Value = 30
$ws.Range("M20").Value = 23
$ws.Range("M21").Value = 23
$ws.Range("M22").Value = 30
$ws.Range("M23").Value = 25
$ws.Range("M24").Value = 28

# Update the per-session "attended/total" label in column H for each
# corresponding block of session rows (110-271), 27 rows per session.
for ($r = 110; $r -le 136; $r++) {
    $ws.Range("H$r").Value = "0/30"
}
for ($r = 137; $r -le 163; $r++) {
    $ws.Range("H$r").Value = "0/23"
}
for ($r = 164; $r -le 190; $r++) {
    $ws.Range("H$r").Value = "0/23"
}
for ($r = 191; $r -le 217; $r++) {
    $ws.Range("H$r").Value = "0/30"
}
for ($r = 218; $r -le 244; $r++) {
    $ws.Range("H$r").Value = "0/25"
}
for ($r = 245; $r -le 271; $r++) {
    $ws.Range("H$r").Value = "0/28"
}
